# Actualización automática 2025-08-27 17:15:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("I11").Value = 813.11
$wsGrupo.Range("M13").Value = 8423.42
$wsGrupo.Range("I23").Value = "2 de 21"

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F11").Value = 2102.43
$wsMensual.Range("F13").Value = 8423.42
$wsMensual.Range("F23").Value = 11734.37

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D8").Value = 855.74
$wsCumplimiento.Range("E8").Value = -230.74
$wsCumplimiento.Range("F8").Value = 1.369184

$wsCumplimiento.Range("D16").Value = 9764.459999999999
$wsCumplimiento.Range("E16").Value = 29012.01
$wsCumplimiento.Range("F16").Value = 0.2518140511500918

$wsCumplimiento.Range("D19").Value = 11734.37
$wsCumplimiento.Range("E19").Value = 47653.85762291768
$wsCumplimiento.Range("F19").Value = 0.197587475997882
